$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit message: "Todos los 0 los cambiamos a 1" -> every numeric 0 in the
# sheet's used range becomes 1 (text cells never equal the number 0, so this
# only ever touches the numeric data column).
$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq 0) {
            $cell.Value = 1
        }
    }
}

# Update the view/selection state as captured in the saved workbook
$ws.Range("E49:E57").Select()
